$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 47, shifting existing rows 47-49 down to 48-50
$ws.Rows.Item(47).Insert()

# Populate the newly inserted row 47 with the new data entry
$ws.Cells.Item(47, 1).Value = 10
$ws.Cells.Item(47, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(47, 3).Value = "La Araucanía"
$ws.Cells.Item(47, 4).Value = 44615
$ws.Cells.Item(47, 5).Value = 9
$ws.Cells.Item(47, 6).Value = 100114002
$ws.Cells.Item(47, 7).Value = "Camote"
$ws.Cells.Item(47, 8).Value = "Sin especificar"
$ws.Cells.Item(47, 9).Value = "Primera"
$ws.Cells.Item(47, 10).Value = 30
$ws.Cells.Item(47, 11).Value = 18000
$ws.Cells.Item(47, 12).Value = 18000
$ws.Cells.Item(47, 13).Value = 18000
$ws.Cells.Item(47, 14).Value = "$/malla 20 kilos"
$ws.Cells.Item(47, 15).Value = "Perú"
$ws.Cells.Item(47, 16).Value = 900
$ws.Cells.Item(47, 17).Value = 20
$ws.Cells.Item(47, 18).Value = "Hortaliza"
